# Update the cryptocurrency price table (rows 2-51) on Sheet1 to reflect
# the latest scrape from coinranking.com. A new coin (OKB) was inserted at
# row 9, shifting Dogecoin..Aave down by one row, and Elrond (previously at
# row 51) fell off the bottom of the table. Price/volume figures were
# refreshed for (almost) every row.
#
# Cells in the Price column (D) that look like plain numbers (e.g. "1.003")
# must still be stored as text, matching the source file's inlineStr cells,
# so we switch those specific cells to a text number format immediately
# before writing them (Excel's COM layer would otherwise silently convert
# a numeric-looking string into a real number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.520.29'
$ws.Cells.Item(3, 4).Value = '1.812.14'
$ws.Cells.Item(3, 5).Value = '  +0.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  -0.49%  '
$ws.Cells.Item(5, 5).Value = '  -0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '305.38'
$ws.Cells.Item(6, 5).Value = '  -0.91%  '
$ws.Cells.Item(7, 5).Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3595'
$ws.Cells.Item(8, 5).Value = '  -1.80%  '
$ws.Cells.Item(9, 2).Value = 'OKB'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '46.29'
$ws.Cells.Item(9, 5).Value = '  +2.54%  '
$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.07106'
$ws.Cells.Item(10, 5).Value = '  -0.28%  '
$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.8932'
$ws.Cells.Item(11, 5).Value = '  +1.63%  '
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.07722'
$ws.Cells.Item(12, 5).Value = '  -0.38%  '
$ws.Cells.Item(13, 2).Value = 'Solana'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '19.28'
$ws.Cells.Item(13, 5).Value = '  -0.38%  '
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value = '1.826.07'
$ws.Cells.Item(14, 5).Value = '  +1.50%  '
$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '5.256'
$ws.Cells.Item(15, 5).Value = '  -0.58%  '
$ws.Cells.Item(16, 2).Value = 'Chainlink'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '6.296'
$ws.Cells.Item(16, 5).Value = '  -1.23%  '
$ws.Cells.Item(17, 2).Value = 'Litecoin'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '85.96'
$ws.Cells.Item(17, 5).Value = '  -0.68%  '
$ws.Cells.Item(18, 2).Value = 'BinanceUSD'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D18').NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '1.006'
$ws.Cells.Item(18, 5).Value = '  -0.39%  '
$ws.Cells.Item(19, 2).Value = 'ShibaInu'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.000008545'
$ws.Cells.Item(19, 5).Value = '  -0.37%  '
$ws.Cells.Item(20, 2).Value = 'Dai'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '1.003'
$ws.Cells.Item(20, 5).Value = '  -0.39%  '
$ws.Cells.Item(21, 2).Value = 'WrappedBTC'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(21, 4).Value = '26.561.93'
$ws.Cells.Item(21, 5).Value = '  +0.01%  '
$ws.Cells.Item(22, 2).Value = 'Avalanche'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '14.13'
$ws.Cells.Item(22, 5).Value = '  -0.69%  '
$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '4.956'
$ws.Cells.Item(23, 5).Value = '  -1.06%  '
$ws.Cells.Item(24, 2).Value = 'Cosmos'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '10.49'
$ws.Cells.Item(24, 5).Value = '  +0.14%  '
$ws.Cells.Item(25, 2).Value = 'Toncoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '1.921'
$ws.Cells.Item(25, 5).Value = '  -3.47%  '
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '152.07'
$ws.Cells.Item(26, 5).Value = '  +0.41%  '
$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '17.78'
$ws.Cells.Item(27, 5).Value = '  -0.85%  '
$ws.Cells.Item(28, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.023'
$ws.Cells.Item(28, 5).Value = '  -1.79%  '
$ws.Cells.Item(29, 2).Value = 'BitcoinCash'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '112.37'
$ws.Cells.Item(29, 5).Value = '  -0.31%  '
$ws.Cells.Item(30, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '4.820'
$ws.Cells.Item(30, 5).Value = '  -0.35%  '
$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.08708'
$ws.Cells.Item(31, 5).Value = '  +0.31%  '
$ws.Cells.Item(32, 2).Value = 'HuobiToken'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.133'
$ws.Cells.Item(32, 5).Value = '  +2.34%  '
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.7375'
$ws.Cells.Item(33, 5).Value = '  +0.81%  '
$ws.Cells.Item(34, 2).Value = 'Filecoin'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.423'
$ws.Cells.Item(34, 5).Value = '  -1.86%  '
$ws.Cells.Item(35, 2).Value = 'RenderToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D35').NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '2.713'
$ws.Cells.Item(35, 5).Value = '  +1.91%  '
$ws.Cells.Item(36, 2).Value = 'ARBITRUM'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.109'
$ws.Cells.Item(36, 5).Value = '  -0.78%  '
$ws.Cells.Item(37, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '1.071'
$ws.Cells.Item(37, 5).Value = '  -1.11%  '
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.01934'
$ws.Cells.Item(38, 5).Value = '  -0.86%  '
$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '2.913'
$ws.Cells.Item(39, 5).Value = '  +0.55%  '
$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.05080'
$ws.Cells.Item(40, 5).Value = '  -0.51%  '
$ws.Cells.Item(41, 2).Value = 'TheSandbox'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.5072'
$ws.Cells.Item(41, 5).Value = '  +1.74%  '
$ws.Cells.Item(42, 2).Value = 'FraxShare'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '6.786'
$ws.Cells.Item(42, 5).Value = '  -2.57%  '
$ws.Cells.Item(43, 2).Value = 'Algorand'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.1507'
$ws.Cells.Item(43, 5).Value = '  -3.40%  '
$ws.Cells.Item(44, 2).Value = 'Aptos'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '8.020'
$ws.Cells.Item(44, 5).Value = '  -1.71%  '
$ws.Cells.Item(45, 2).Value = 'Decentraland'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.4661'
$ws.Cells.Item(45, 5).Value = '  +1.40%  '
$ws.Cells.Item(46, 2).Value = 'PaxDollar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '1.003'
$ws.Cells.Item(46, 5).Value = '  -0.44%  '
$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '9.933'
$ws.Cells.Item(47, 5).Value = '  -0.07%  '
$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '98.99'
$ws.Cells.Item(48, 5).Value = '  -1.97%  '
$ws.Cells.Item(49, 2).Value = 'NEARProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.564'
$ws.Cells.Item(49, 5).Value = '  -1.55%  '
$ws.Cells.Item(50, 2).Value = 'Cronos'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.05993'
$ws.Cells.Item(50, 5).Value = '  -0.12%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '63.70'
$ws.Cells.Item(51, 5).Value = '  -0.95%  '
